$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "C:\Users\BraxtonWright\Desktop\Ark files"
$ws.Range("B2").Value = "H:\my games\Steam\Game notes\Ark files"
$ws.Range("C2").Value = "Ark notes"

$ws.Range("A3").Value = "C:\Users\BraxtonWright\Desktop"
$ws.Range("C3").Value = "Satisfactory notes"
$ws.Range("B3").Value = "H:\my games\Steam\Game notes\Satisfactory files"
$ws.Range("D3").Value = "^Satisfactory.txt$"

$ws.Range("D4").Select()
